$d = $word.ActiveDocument

# --- Move the "_GoBack" bookmark from the start of the document to the
# --- end of the "Work Experience" heading paragraph (right after the
# --- run containing "xperience", immediately before that paragraph's
# --- mark).
#
# Word (and this emulation) only ever keeps one bookmark per name, so
# calling Bookmarks.Add("_GoBack", <new range>) relocates the existing
# "_GoBack" bookmark rather than creating a duplicate - exactly like
# Word silently re-dropping the cursor-position bookmark when the
# document is saved after an edit elsewhere.
#
# Bookmarks.Add() on a truly collapsed range sitting exactly at
# "end-of-paragraph-text, just before the pilcrow" does not relocate the
# bookmark reliably. To sidestep that, temporarily insert a one
# character placeholder at that exact spot, wrap a (non-collapsed)
# bookmark range around the placeholder, then delete the placeholder
# again - the bookmark collapses back down and stays put.

$rng = $d.Content
$found = $rng.Find.Execute("Work Experience", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Work Experience' heading paragraph"
}
$rng.Collapse(0)  # wdCollapseEnd -> right before the paragraph mark
$insertPos = $rng.Start

$marker = $d.Range($insertPos, $insertPos)
$marker.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($insertPos, $insertPos + 1)
$delRange.Delete()
